$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set NumberFormat to Text for price cells whose values are plain numeric-looking
# strings, so Excel preserves the exact original text (trailing zeros, etc.)
# instead of silently re-parsing them as numbers.
$textCells = @("D4","D6","D7","D8","D9","D10","D11","D12","D13","D15","D18","D19","D20","D21","D22","D24","D26","D27","D28","D29","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated price (column D) values
$ws.Range("D4").Value = '0.9998'
$ws.Range("D6").Value = '289.42'
$ws.Range("D7").Value = '0.3936'
$ws.Range("D8").Value = '0.3174'
$ws.Range("D9").Value = '42.23'
$ws.Range("D10").Value = '0.07175'
$ws.Range("D11").Value = '1.056'
$ws.Range("D12").Value = '1.000'
$ws.Range("D13").Value = '5.685'
$ws.Range("D15").Value = '6.597'
$ws.Range("D18").Value = '0.06606'
$ws.Range("D19").Value = '83.72'
$ws.Range("D20").Value = '0.9997'
$ws.Range("D21").Value = '6.120'
$ws.Range("D22").Value = '15.46'
$ws.Range("D24").Value = '2.350'
$ws.Range("D26").Value = '2.352'
$ws.Range("D27").Value = '149.68'
$ws.Range("D28").Value = '18.33'
$ws.Range("D29").Value = '4.842'
$ws.Range("D31").Value = '117.00'
$ws.Range("D32").Value = '6.074'
$ws.Range("D33").Value = '0.9401'
$ws.Range("D34").Value = '0.08129'
$ws.Range("D35").Value = '8.521'
$ws.Range("D36").Value = '5.168'
$ws.Range("D37").Value = '0.06016'
$ws.Range("D38").Value = '0.02227'
$ws.Range("D39").Value = '1.455'
$ws.Range("D40").Value = '0.2027'
$ws.Range("D42").Value = '10.95'
$ws.Range("D43").Value = '0.9994'
$ws.Range("D44").Value = '0.5774'
$ws.Range("D45").Value = '13.01'
$ws.Range("D46").Value = '3.714'
$ws.Range("D47").Value = '0.5515'
$ws.Range("D48").Value = '1.166'
$ws.Range("D49").Value = '1.880'
$ws.Range("D50").Value = '115.98'
$ws.Range("D51").Value = '0.06690'

$ws.Range("D2").Value = '21.663.15'
$ws.Range("D3").Value = '1.534.78'
$ws.Range("D16").Value = '1.545.76'
$ws.Range("D25").Value = '21.669.74'
$ws.Range("D30").Value = '1.721.17'

# Updated volume/percentage change (column E) values
$ws.Range("E2").Value = '  -1.87%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("E6").Value = '  +0.63%  '
$ws.Range("E7").Value = '  +1.79%  '
$ws.Range("E8").Value = '  -2.13%  '
$ws.Range("E9").Value = '  -1.12%  '
$ws.Range("E10").Value = '  -2.62%  '
$ws.Range("E11").Value = '  -5.88%  '
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("E13").Value = '  -0.23%  '
$ws.Range("E14").Value = '  -4.06%  '
$ws.Range("E15").Value = '  -3.04%  '
$ws.Range("E16").Value = '  -0.62%  '
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("E19").Value = '  -1.85%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  -4.49%  '
$ws.Range("E22").Value = '  -3.46%  '
$ws.Range("E23").Value = '  -6.89%  '
$ws.Range("E25").Value = '  -1.89%  '
$ws.Range("E26").Value = '  -7.87%  '
$ws.Range("E27").Value = '  +0.11%  '
$ws.Range("E28").Value = '  -2.95%  '
$ws.Range("E29").Value = '  -0.47%  '
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("E31").Value = '  -3.24%  '
$ws.Range("E32").Value = '  +4.09%  '
$ws.Range("E33").Value = '  -15.18%  '
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("E35").Value = '  -8.96%  '
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("E37").Value = '  -4.19%  '
$ws.Range("E38").Value = '  -3.34%  '
$ws.Range("E39").Value = '  -14.46%  '
$ws.Range("E40").Value = '  -4.12%  '
$ws.Range("E41").Value = '  -3.75%  '
$ws.Range("E42").Value = '  +0.67%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E44").Value = '  -3.15%  '
$ws.Range("E45").Value = '  -3.53%  '
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("E47").Value = '  -4.21%  '
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("E49").Value = '  -2.73%  '
$ws.Range("E50").Value = '  -2.59%  '
$ws.Range("E51").Value = '  -2.92%  '
